# feat: add 2022-Q3 data
#
# Target state:
#  - "总计" (summary) sheet: row 2 becomes the new 2022-Q3 totals, and the
#    previous 2021-Q1 totals are pushed down into a new row 3.
#  - A brand-new "2022-Q3" worksheet holding that quarter's fund-holdings
#    table is inserted right after "总计" (i.e. before "2021-Q1").
#
# Implementation note: rather than inserting a blank sheet for "2022-Q3",
# we duplicate the existing "2021-Q1" sheet (so the untouched "2021-Q1"
# copy keeps its original look/formatting) and turn the original sheet
# into "2022-Q3" by overwriting its content — this mirrors how the
# change was actually authored (the original sheet's underlying part is
# reused for the new quarter, and the duplicate becomes the "new" copy
# of "2021-Q1").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push the current row 2 (2021-Q1) data down into a new row 3, matching
# the formatting already used on row 2's index cell (A2).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q1"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

# Overwrite row 2 with the new 2022-Q3 totals
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("D2").Value = 0.03

# ---------------------------------------------------------------------
# 2. Duplicate "2021-Q1" so the copy keeps that quarter's data untouched,
#    then repurpose the original sheet as "2022-Q3" with new content.
# ---------------------------------------------------------------------
$original2021 = $wb.Worksheets.Item("2021-Q1")
$original2021.Copy($null, $original2021)
$newCopy = $wb.Worksheets.Item("2021-Q1 (2)")

$newSheet = $original2021
$newSheet.Name = "2022-Q3"
$newCopy.Name = "2021-Q1"

# ---- Header row (plain text, no auto-conversion risk) ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- Data row 2 ----
$newSheet.Range("A2").Value = 0
$newSheet.Range("C2").Value = "国海量化优选一年持有股票C"
$newSheet.Range("H2").Value = 4

# ---- Data row 3 ----
$newSheet.Range("A3").Value = 1
$newSheet.Range("C3").Value = "国海量化优选一年持有股票A"
$newSheet.Range("H3").Value = 4

# ---- Numeric-looking identifiers/figures that must stay TEXT (they are
#      strings in the source data: fund codes and percentages). Force the
#      cell format to Text before writing so Excel doesn't coerce them to
#      numbers, then drop back to the default "Normal" cell style so no
#      stray formatting is left behind on these cells. ----
$plainTextCells = @{
    "B2" = "970042"
    "D2" = "7.16"
    "E2" = "87.31"
    "F2" = "0.34"
    "G2" = "0.0243"
    "B3" = "970041"
    "D3" = "0.63"
    "E3" = "87.31"
    "F3" = "0.34"
    "G3" = "0.0021"
}
foreach ($addr in $plainTextCells.Keys) {
    $newSheet.Range($addr).NumberFormat = "@"
    $newSheet.Range($addr).Value = $plainTextCells[$addr]
    $newSheet.Range($addr).Style = "Normal"
}

# ---- Formatting: reuse the same style already used on 总计!A2 for the
#      header row and the index column, matching the source workbook
#      (no new styles are introduced). ----
$totalSheet.Range("A2").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A3").PasteSpecial(-4122)
